$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Clear cells in rows 3-5 that were removed (G:AG range)
$ws.Range("AC3").ClearContents()
$ws.Range("AG3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("Y3").ClearContents()

$ws.Range("AC4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("P4").ClearContents()
$ws.Range("Q4").ClearContents()
$ws.Range("R4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("W4").ClearContents()
$ws.Range("Y4").ClearContents()

$ws.Range("AC5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("Y5").ClearContents()

# Set new values for rows 6-12 (mean, std, min, 25%, 50%, 75%, max)
$ws.Range("G6").Value = 48.3436898471683
$ws.Range("H6").Value = 42.93699873937364
$ws.Range("I6").Value = 48.49055426504608
$ws.Range("J6").Value = 52.70447787809864
$ws.Range("K6").Value = 2407.262708239891
$ws.Range("L6").Value = 2426.960880160337
$ws.Range("M6").Value = 2717.543403275555
$ws.Range("N6").Value = 154397.1234062055
$ws.Range("O6").Value = 5.33389643857483
$ws.Range("P6").Value = 17.37412241208274
$ws.Range("Q6").Value = 82.19905607144214
$ws.Range("R6").Value = 0.8584311403732222
$ws.Range("S6").Value = 0.661922958745513
$ws.Range("T6").Value = 5.972326062666207
$ws.Range("U6").Value = 26.18777622851844
$ws.Range("V6").Value = 6908805137.07477
$ws.Range("W6").Value = 27.51716344048013
$ws.Range("Y6").Value = 0.08024453236167123
$ws.Range("AC6").Value = 2.767140727173219
$ws.Range("AG6").Value = 22.2725576081287

$ws.Range("G7").Value = 17.16147464662168
$ws.Range("H7").Value = 23.70901205442084
$ws.Range("I7").Value = 18.91038928922768
$ws.Range("J7").Value = 22.69754271402129
$ws.Range("K7").Value = 1910.874363925357
$ws.Range("L7").Value = 1769.776694781962
$ws.Range("M7").Value = 1684.119621374067
$ws.Range("N7").Value = 129965.4632746056
$ws.Range("O7").Value = 9.943481437683465
$ws.Range("P7").Value = 35.19920109302316
$ws.Range("Q7").Value = 859.3073582744092
$ws.Range("R7").Value = 0.3777431378329725
$ws.Range("S7").Value = 0.305782261164173
$ws.Range("T7").Value = 8.734407333577249
$ws.Range("U7").Value = 19.36566613006569
$ws.Range("V7").Value = 11761719015.21445
$ws.Range("W7").Value = 25.77032520041962
$ws.Range("Y7").Value = 0.06852723214106546
$ws.Range("AC7").Value = 3.073625885602018
$ws.Range("AG7").Value = 73.63039604550032

$ws.Range("G8").Value = 10.3020035377314
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 4.96888175065247
$ws.Range("J8").Value = 2.04954954954955
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 14.88426290172922
$ws.Range("N8").Value = 1093.364789999466
$ws.Range("O8").Value = 0.05
$ws.Range("P8").Value = 0.745169718272074
$ws.Range("Q8").Value = 1.86292715346454
$ws.Range("R8").Value = -0.06247586302255
$ws.Range("S8").Value = 0.001413338582652
$ws.Range("T8").Value = -15.065971452
$ws.Range("U8").Value = 8.71281742681019
$ws.Range("V8").Value = 141908221.736258
$ws.Range("W8").Value = 8.68223590341081
$ws.Range("Y8").Value = [double]"7.078209147174905e-05"
$ws.Range("AC8").Value = 0.031542391616566
$ws.Range("AG8").Value = 0.06676437641096522

$ws.Range("G9").Value = 35.31561971254833
$ws.Range("H9").Value = 23.03394184904595
$ws.Range("I9").Value = 34.52821320855253
$ws.Range("J9").Value = 35.28173510278123
$ws.Range("K9").Value = 853.3856318232822
$ws.Range("L9").Value = 890.1100922164812
$ws.Range("M9").Value = 1369.587820076888
$ws.Range("N9").Value = 44046.65264847894
$ws.Range("O9").Value = 0.60625
$ws.Range("P9").Value = 5.119856057319669
$ws.Range("Q9").Value = 10.560535286
$ws.Range("R9").Value = 0.6398918030292789
$ws.Range("S9").Value = 0.452307957022013
$ws.Range("T9").Value = 1.481610028775
$ws.Range("U9").Value = 16.7236466118598
$ws.Range("V9").Value = 888541699.6894521
$ws.Range("W9").Value = 18.2679219130421
$ws.Range("Y9").Value = 0.04412667376794756
$ws.Range("AC9").Value = 0.8053139934798538
$ws.Range("AG9").Value = 2.623179182637529

$ws.Range("G10").Value = 48.641297233402
$ws.Range("H10").Value = 40.87287419410555
$ws.Range("I10").Value = 48.3022370800429
$ws.Range("J10").Value = 54.7248857801159
$ws.Range("K10").Value = 1850.039300710159
$ws.Range("L10").Value = 2218.393380161163
$ws.Range("M10").Value = 2569.4559280577
$ws.Range("N10").Value = 115084.1789179667
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 10.1644475225559
$ws.Range("Q10").Value = 15.6640487231732
$ws.Range("R10").Value = 0.871722244882923
$ws.Range("S10").Value = 0.6689679885847239
$ws.Range("T10").Value = 4.271491717
$ws.Range("U10").Value = 20.3795031816288
$ws.Range("V10").Value = 2006301764.7774
$ws.Range("W10").Value = 21.657082377917
$ws.Range("Y10").Value = 0.06384394972145617
$ws.Range("AC10").Value = 1.647846535
$ws.Range("AG10").Value = 7.044884615320818

$ws.Range("G11").Value = 62.51485132655755
$ws.Range("H11").Value = 62.53446260842355
$ws.Range("I11").Value = 62.7283928166508
$ws.Range("J11").Value = 68.23965141612197
$ws.Range("K11").Value = 3885.23083918973
$ws.Range("L11").Value = 3613.386255954538
$ws.Range("M11").Value = 3817.096439811206
$ws.Range("N11").Value = 244315.3523009635
$ws.Range("O11").Value = 3.9275
$ws.Range("P11").Value = 17.61055795925
$ws.Range("Q11").Value = 22.66210191
$ws.Range("R11").Value = 1.13712072991543
$ws.Range("S11").Value = 0.860110120800776
$ws.Range("T11").Value = 7.3692958015
$ws.Range("U11").Value = 31.4613034542723
$ws.Range("V11").Value = 7159835140.950863
$ws.Range("W11").Value = 32.1906132064199
$ws.Range("Y11").Value = 0.09469218619755916
$ws.Range("AC11").Value = 3.53125050266558
$ws.Range("AG11").Value = 16.42101713039693

$ws.Range("G12").Value = 82.0670604173895
$ws.Range("H12").Value = 95.1104018863559
$ws.Range("I12").Value = 84.35048299428
$ws.Range("J12").Value = 97.7378006872852
$ws.Range("K12").Value = 7544.707386704818
$ws.Range("L12").Value = 7820.196602272073
$ws.Range("M12").Value = 7944.860832502553
$ws.Range("N12").Value = 552721.8493276552
$ws.Range("O12").Value = 69.58
$ws.Range("P12").Value = 417.816983778207
$ws.Range("Q12").Value = 14127.867363161
$ws.Range("R12").Value = 2.04619443595222
$ws.Range("S12").Value = 1.42094148620397
$ws.Range("T12").Value = 72.996943993
$ws.Range("U12").Value = 151.057635055307
$ws.Range("V12").Value = 65335707786.2141
$ws.Range("W12").Value = 207.701129958666
$ws.Range("Y12").Value = 0.5367896421179276
$ws.Range("AC12").Value = 26.1724181013647
$ws.Range("AG12").Value = 753.928225797619
